# Adds newly-researched cabinet "Position" data (column J) for the
# Conservative and Labour front-benches, fills in Rishi Sunak's "Power"
# score, and adds a full profile (Introduction / Random fact / Social
# media) for Liz Truss, including a real hyperlink with rich text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Labour shadow cabinet (rows 17-22, column J) -------------------------
$ws.Range("J17").Value = "Shadow Secretary of State for Health and Social Care"
$ws.Range("J18").Value = "Shadow Secretary of State for Energy Security and Net Zero"
$ws.Range("J19").Value = "Shadow Secretary of State for Foreign, Commonwealth and Development Affairs"
$ws.Range("J20").Value = "Shadow Chancellor of the Duchy of Lancaster and National Campaign Coordinator"
$ws.Range("J21").Value = "Shadow Minister without Portfolio (Cabinet Office)"
$ws.Range("J22").Value = "Shadow Paymaster General (Cabinet Office)"

# --- Conservative cabinet (rows 2-11, column J) ---------------------------
$ws.Range("J3").Value = "Prime Minister"
$ws.Range("J7").Value = "Secretary of State for the Home Department"
$ws.Range("J9").Value = "Secretary of State for Defence"
$ws.Range("J10").Value = "Lord Chancellor and Secretary of State for Justice"
$ws.Range("J11").Value = "Secretary of State for Science, Innovation and Technology"
$ws.Range("J8").Value = "Chancellor of the Duchy of Lancaster, and Secretary of State in the Cabinet Office; Deputy Prime Minister"
$ws.Range("J6").Value = "Secretary of State for Foreign, Commonwealth and Development Affairs"
$ws.Range("J4").Value = "Ex-Prime Minister"
$ws.Range("J2").Value = "Former Prime Minister"

# Rishi Sunak's "Power" score
$ws.Range("B3").Value = 100

# --- Liz Truss profile (row 4) --------------------------------------------
$ws.Range("D4").Value = "Elizabeth Truss was Prime Minister from 6 September 2022 to 25 October 2022. She was previously Secretary of State for Foreign, Commonwealth and Development Affairs from 15 September 2021. She was appointed Minister for Women and Equalities on 10 September 2019. She was elected as the Conservative MP for south west Norfolk in 2010."

$ws.Range("E4").Value = "Elizabeth was previously Deputy Director at Reform. She also worked in the energy and telecommunications industry for 10 years as a commercial manager and economics director, and is a qualified management accountant."
$ws.Range("E4").Font.Name = "Arial"
$ws.Range("E4").Font.Color = 789515

# The real hyperlink (creates the workbook's "Hyperlink" cell style) plus
# the two-line rich-text display (twitter handle on line 1, bold LinkedIn
# URL on line 2). The "display" argument both sets the hyperlink's stored
# display text AND the cell's value, so it must be passed here rather than
# via .Value first.
$display = "https://twitter.com/trussliz`nhttps://www.linkedin.com/in/liz-truss/"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://twitter.com/trussliz", $null, $null, $display)

# J5 happens to carry the (unused) Hyperlink cell style in the source file;
# grab the plain "Hyperlink" style right away, before F4 below gets its own
# (wrap-text) variant of that same style.
$ws.Range("J5").Value = "Chancellor of the Exchequer"
$ws.Range("J5").Style = "Hyperlink"

$chars = $ws.Range("F4").Characters(30, 38)
$chars.Font.Bold = $true
$chars.Font.Underline = $true
$chars.Font.Color = 12673797
$ws.Range("F4").WrapText = $true

$ws.Rows(4).RowHeight = 119

# --- Column widths / selection --------------------------------------------
$ws.Columns("D").ColumnWidth = 152.5
$ws.Columns("J").ColumnWidth = 102.5

[void]$ws.Range("D7").Select()
